$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.408.44"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "2.190.85"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  +0.00%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "255.30"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +4.86%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.629"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +0.70%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "68.59"
$cell.ClearFormats()
$ws.Range("E7").Value = "  -2.49%  "
$ws.Range("E8").Value = "  +0.06%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.577"
$cell.ClearFormats()
$ws.Range("E9").Value = "  +5.52%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "38.01"
$cell.ClearFormats()
$ws.Range("E10").Value = "  +2.86%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "59.06"
$cell.ClearFormats()
$ws.Range("E11").Value = "  +2.51%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.0939"
$cell.ClearFormats()
$ws.Range("E12").Value = "  -1.14%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "7.16"
$cell.ClearFormats()
$ws.Range("E13").Value = "  +7.80%  "
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "2.516.78"
$ws.Range("E15").Value = "  -0.54%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.878"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +4.93%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "14.54"
$cell.ClearFormats()
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("D18").Value = "2.203.25"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "41.322.74"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").Value = "0.0₃0958"
$ws.Range("E20").Value = "  +1.41%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "6.19"
$cell.ClearFormats()
$ws.Range("E21").Value = "  +1.98%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "72.09"
$cell.ClearFormats()
$ws.Range("E22").Value = "  -0.74%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "232.61"
$cell.ClearFormats()
$ws.Range("E23").Value = "  +0.86%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.02"
$cell.ClearFormats()
$ws.Range("E24").Value = "  +0.72%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "3.94"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +10.19%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "11.78"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +21.19%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").Value = "  +5.41%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "3.73"
$cell.ClearFormats()
$ws.Range("E29").Value = "  -4.47%  "
$ws.Range("E30").Value = "  -0.41%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "168.95"
$cell.ClearFormats()
$ws.Range("E31").Value = "  -1.16%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "20.70"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +1.91%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.118"
$cell.ClearFormats()
$ws.Range("E33").Value = "  -0.51%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.0757"
$cell.ClearFormats()
$ws.Range("E34").Value = "  +7.50%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.124"
$cell.ClearFormats()
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  +6.29%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "26.67"
$cell.ClearFormats()
$ws.Range("E37").Value = "  +12.41%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "4.64"
$cell.ClearFormats()
$ws.Range("E38").Value = "  +0.97%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "4.16"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +8.00%  "
$ws.Range("E40").Value = "  +9.74%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "12.65"
$cell.ClearFormats()
$ws.Range("E41").Value = "  +19.41%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "2.21"
$cell.ClearFormats()
$ws.Range("E42").Value = "  -2.61%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "5.69"
$cell.ClearFormats()
$ws.Range("E43").Value = "  -1.67%  "
$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "64.47"
$cell.ClearFormats()
$ws.Range("E44").Value = "  +1.87%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "5.12"
$cell.ClearFormats()
$ws.Range("E45").Value = "  +5.09%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.202"
$cell.ClearFormats()
$ws.Range("E46").Value = "  +3.60%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "8.68"
$cell.ClearFormats()
$ws.Range("E47").Value = "  +1.16%  "
$ws.Range("E48").Value = "  +2.37%  "
$ws.Range("E49").Value = "  +0.40%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.15"
$cell.ClearFormats()
$ws.Range("E50").Value = "  +5.69%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.18"
$cell.ClearFormats()
$ws.Range("E51").Value = "  +1.10%  "
